$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.536.19'
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.598.99'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.36'
$ws.Range("E5").Value = '  +0.56%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.45'
$ws.Range("E6").Value = '  +3.33%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '8.05'
$ws.Range("E9").Value = '  +1.69%  '
$ws.Range("E10").Value = '  -0.19%  '
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.210.25'
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("E13").Value = '  +0.74%  '
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.598.05'
$ws.Range("E15").Value = '  +0.94%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.602.67'
$ws.Range("E16").Value = '  +0.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.117'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.58'
$ws.Range("E18").Value = '  +1.92%  '
$ws.Range("E19").Value = '  +3.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.12'
$ws.Range("E20").Value = '  +1.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '427.43'
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.620'
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '78.76'
$ws.Range("E23").Value = '  -0.33%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("E25").Value = '  +2.50%  '
$ws.Range("E26").Value = '  +4.67%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.43'
$ws.Range("E27").Value = '  +3.57%  '
$ws.Range("E28").Value = '  +0.73%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.596.97'
$ws.Range("E30").Value = '  +1.12%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.46'
$ws.Range("E32").Value = '  -0.51%  '
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.157'
$ws.Range("E33").Value = '  +2.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.65'
$ws.Range("E36").Value = '  +0.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '177.44'
$ws.Range("E39").Value = '  +1.05%  '
$ws.Range("E40").Value = '  +0.86%  '
$ws.Range("E42").Value = '  -2.16%  '
$ws.Range("E43").Value = '  +9.87%  '
$ws.Range("E45").Value = '  -1.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.05'
$ws.Range("E46").Value = '  -2.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.07'
$ws.Range("E47").Value = '  +2.04%  '
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.954'
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.429.70'
$ws.Range("E50").Value = '  +5.51%  '
